# Rename the three embedded logo pictures (the wp:docPr / pic:cNvPr "name"
# attribute only -- ids, descr, and the underlying media parts are untouched):
#   footer1.xml (id=3, PearsonLogo) : image2.png -> image1.png
#   footer2.xml (id=2, PearsonLogo) : image2.png -> image1.png
#   header1.xml (id=1, BTec_Logo)   : image1.jpg -> image2.jpg
#
# InlineShape.Name is not writable through this object model (and
# ConvertToShape()/ConvertToInlineShape() re-homes the drawing as a floating
# anchored shape, which would corrupt the surrounding wp:inline markup), so
# the rename is done by round-tripping the package through
# Document.WordOpenXML and doing a targeted text substitution on the
# wp:docPr/pic:cNvPr name="..." attributes.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image2.png"', 'name="image1.png"')
$xml = $xml.Replace('name="image1.jpg"', 'name="image2.jpg"')

$d.WordOpenXML = $xml

Write-Output "Renamed inline picture name attributes"
